$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "AB Bank" sheet data edits
# ------------------------------------------------------------------
$abBank = $wb.Worksheets.Item("AB Bank")

# B3: 7581 -> 7582 (keeps its existing #,##0 style)
$abBank.Range("B3").Value = 7582

$cols = @("B","C","D","E","F","G","H","I","J","K")

# Row 4 (Shareholders' equity): fill with 0, plain (no number format)
foreach ($col in $cols) {
    $abBank.Range($col + "4").Value = 0
}

# Row 9 (Earning assets): fill with 0, #,##0 style
foreach ($col in $cols) {
    $abBank.Range($col + "9").NumberFormat = "#,##0"
    $abBank.Range($col + "9").Value = 0
}

# Row 10 (Total assets): specific values, #,##0 style
$row10 = @{ "B" = 322526; "C" = 314565; "D" = 314836; "E" = 285010; "F" = 246331; "G" = 208006; "H" = 173842; "I" = 152963; "J" = 132691; "K" = 106912 }
foreach ($col in $cols) {
    $abBank.Range($col + "10").NumberFormat = "#,##0"
    $abBank.Range($col + "10").Value = $row10[$col]
}

# Row 11 (Total liabilities): fill with 0, #,##0 style
foreach ($col in $cols) {
    $abBank.Range($col + "11").NumberFormat = "#,##0"
    $abBank.Range($col + "11").Value = 0
}

# Row 15 (Guarantee Business): fill with 0, #,##0 style
foreach ($col in $cols) {
    $abBank.Range($col + "15").NumberFormat = "#,##0"
    $abBank.Range($col + "15").Value = 0
}

# ------------------------------------------------------------------
# "MTB" sheet: move the (empty, red-font styled) cell from B16 to B15
# ------------------------------------------------------------------
$mtb = $wb.Worksheets.Item("MTB")
$mtb.Range("B15").Font.Color = 255
$mtb.Range("B16").Clear()

# ------------------------------------------------------------------
# Sheet selections / active sheet.
# Selecting a range switches to that worksheet, so the last sheet
# touched becomes the active / tabSelected one. Do "AB Bank" last so
# it ends up active (matches activeTab="4" / tabSelected="1").
# ------------------------------------------------------------------
$primeBank = $wb.Worksheets.Item("Prime Bank")
$null = $primeBank.Range("B9").Select()

$dhakaBank = $wb.Worksheets.Item("Dhaka Bank")
$null = $dhakaBank.Range("B17").Select()

$null = $mtb.Range("B17").Select()

$null = $abBank.Range("M13").Select()
